$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("_reportFile", "blocks", "path", "status", "template", "timestamp", "title", "url"),
  @("index.report.json", "[`"hero-banner`",`"cards-promo`",`"cards-promo`",`"cards-value`",`"carousel-industry`",`"cards-story`",`"accordion-links`"]", "index", "success", "homepage", "2026-02-26T16:46:19.020Z", "AT&T Business – Fiber Internet, Wireless, Phone, IoT, 5G Solutions", "https://www.business.att.com/"),
  @("learn.report.json", "[`"hero-banner`",`"cards-story`",`"cards-story`"]", "learn", "success", "hub-page", "2026-03-02T10:02:52.172Z", "Welcome To Business Insights and Technology Advice", "https://www.business.att.com/learn.html"),
  @("small-business.report.json", "[]", "small-business", "success", "industry-page", "2026-03-02T10:02:14.236Z", "Explore AT&T Internet Devices for Home and Small Business", "https://www.business.att.com/small-business.html"),
  @("support.report.json", "[]", "support", "success", "hub-page", "2026-03-02T10:03:04.405Z", "AT&T Business Support", "https://www.business.att.com/support.html"),
  @("about/why-att-business.report.json", "[`"hero-banner`",`"cards-value`",`"columns-offer`",`"accordion-links`"]", "about/why-att-business", "success", "landing-page", "2026-03-02T10:02:29.654Z", "Why Customers Choose AT&T Business For Expert Guidance and Advice", "https://www.business.att.com/about/why-att-business.html"),
  @("industries/finance.report.json", "[`"hero-banner`",`"hero-banner`"]", "industries/finance", "success", "industry-page", "2026-03-02T10:01:49.111Z", "Financial Technology & Security Solutions at AT&T Business.", "https://www.business.att.com/industries/finance.html"),
  @("industries/healthcare.report.json", "[`"hero-banner`",`"cards-value`",`"columns-offer`"]", "industries/healthcare", "success", "industry-page", "2026-03-02T10:01:33.523Z", "Cybersecurity, Networking & Iot Healthcare Solutions at AT&T.", "https://www.business.att.com/industries/healthcare.html"),
  @("industries/hospitality.report.json", "[`"hero-banner`",`"hero-banner`",`"columns-offer`"]", "industries/hospitality", "success", "industry-page", "2026-03-02T10:01:43.957Z", "Hospitality Technology & Guest Service Solutions at AT&T Business.", "https://www.business.att.com/industries/hospitality.html"),
  @("industries/manufacturing.report.json", "[`"hero-banner`",`"cards-value`",`"columns-offer`",`"columns-offer`",`"columns-offer`"]", "industries/manufacturing", "success", "industry-page", "2026-03-02T10:02:00.717Z", "IOT Network and Technology Solutions for Manufactures at AT&T Business", "https://www.business.att.com/industries/manufacturing.html"),
  @("industries/public-sector.report.json", "[`"hero-banner`",`"hero-banner`",`"hero-banner`",`"hero-banner`",`"columns-offer`",`"columns-offer`"]", "industries/public-sector", "success", "industry-page", "2026-03-02T10:02:08.188Z", "AT&T Public Sector Solutions deliver a faster, smarter, safer future", "https://www.business.att.com/industries/public-sector.html"),
  @("industries/retail.report.json", "[`"hero-banner`",`"cards-value`",`"columns-offer`"]", "industries/retail", "success", "industry-page", "2026-03-02T10:01:39.311Z", "Connected Consumer & Retail Technology Solutions at AT&T Business.", "https://www.business.att.com/industries/retail.html"),
  @("industries/transportation.report.json", "[`"hero-banner`",`"cards-value`",`"columns-offer`"]", "industries/transportation", "success", "industry-page", "2026-03-02T10:01:53.480Z", "Fiber, 5G and other digital solutions for the transportation industry", "https://www.business.att.com/industries/transportation.html"),
  @("learn/customer-stories.report.json", "[`"hero-banner`"]", "learn/customer-stories", "success", "hub-page", "2026-03-02T10:02:58.245Z", "Customer Stories | AT&T Business", "https://www.business.att.com/learn/customer-stories.html"),
  @("offers/switch-to-att-business.report.json", "[`"hero-banner`",`"cards-value`",`"cards-value`",`"cards-product`",`"columns-offer`"]", "offers/switch-to-att-business", "success", "landing-page", "2026-03-02T10:02:36.674Z", "Change your wireless carrier and internet. Switch to AT&T Business", "https://www.business.att.com/offers/switch-to-att-business.html"),
  @("portfolios/business-internet.report.json", "[`"cards-product`",`"cards-value`",`"hero-dark`",`"accordion-links`"]", "portfolios/business-internet", "success", "portfolio-page", "2026-03-02T10:00:19.677Z", "Ultra-Fast Business Internet & Wi-Fi Services | AT&T Business", "https://www.business.att.com/portfolios/business-internet.html"),
  @("portfolios/collaboration.report.json", "[`"hero-banner`",`"cards-promo`",`"cards-product`",`"cards-product`",`"cards-value`"]", "portfolios/collaboration", "success", "portfolio-page", "2026-03-02T10:00:25.797Z", "VOIP, Voice, Video & Chat Solutions at AT&T Business.", "https://www.business.att.com/portfolios/collaboration.html"),
  @("portfolios/cybersecurity.report.json", "[`"hero-banner`",`"cards-product`",`"cards-product`",`"cards-value`",`"columns-offer`",`"columns-offer`"]", "portfolios/cybersecurity", "success", "portfolio-page", "2026-03-02T10:00:32.910Z", "Secure Your Business with Cybersecurity Services | AT&T Business", "https://www.business.att.com/portfolios/cybersecurity.html"),
  @("portfolios/internet-of-things.report.json", "[`"hero-banner`",`"cards-product`",`"cards-product`",`"cards-product`",`"cards-product`",`"cards-value`",`"columns-offer`",`"columns-offer`",`"columns-offer`",`"accordion-links`"]", "portfolios/internet-of-things", "success", "portfolio-page", "2026-03-02T10:00:44.466Z", "IoT Solutions and Platform as a Service at AT&T Business", "https://www.business.att.com/portfolios/internet-of-things.html"),
  @("portfolios/mobility.report.json", "[`"hero-banner`",`"cards-promo`",`"cards-product`",`"cards-value`",`"columns-offer`",`"accordion-links`"]", "portfolios/mobility", "success", "portfolio-page", "2026-03-02T10:00:11.024Z", "AT&T Wireless Solutions, Mobile Business Services & Support", "https://www.business.att.com/portfolios/mobility.html"),
  @("portfolios/networking.report.json", "[`"cards-product`",`"cards-product`",`"cards-product`",`"cards-value`",`"hero-dark`",`"accordion-links`"]", "portfolios/networking", "success", "portfolio-page", "2026-03-02T10:00:38.966Z", "Business Networking Services & Solutions | AT&T Business", "https://www.business.att.com/portfolios/networking.html"),
  @("products/att-dynamic-defense.report.json", "[`"hero-banner`",`"cards-product`",`"columns-offer`",`"accordion-faq`"]", "products/att-dynamic-defense", "success", "product-page", "2026-03-02T10:01:11.789Z", "Enterprise Network Security and Threat Protection", "https://www.business.att.com/products/att-dynamic-defense.html"),
  @("products/business-fiber-internet.report.json", "[`"hero-dark`",`"columns-offer`",`"accordion-faq`",`"accordion-faq`"]", "products/business-fiber-internet", "success", "product-page", "2026-03-02T10:01:06.605Z", "AT&T Business Fiber®: High Speed Business Fiber Internet Service", "https://www.business.att.com/products/business-fiber-internet.html"),
  @("products/sd-wan.report.json", "[`"cards-value`",`"columns-offer`",`"columns-offer`",`"accordion-faq`"]", "products/sd-wan", "success", "product-page", "2026-03-02T10:01:17.688Z", "SD-WAN Service Provider - Software Defined Wide Area Network at AT&T Business.", "https://www.business.att.com/products/sd-wan.html"),
  @("products/wireless-plans.report.json", "[`"hero-banner`",`"cards-product`",`"cards-product`",`"cards-product`",`"cards-value`",`"cards-value`",`"hero-dark`",`"hero-dark`",`"hero-dark`",`"columns-offer`",`"accordion-faq`"]", "products/wireless-plans", "success", "product-page", "2026-03-02T13:40:12.511Z", "Business Cell Phone Plans | AT&T Business Wireless & Mobile", "https://www.business.att.com/products/wireless-plans.html"),
  @("learn/customer-stories/parkville.report.json", "[`"cards-value`",`"cards-promo`",`"columns-offer`"]", "learn/customer-stories/parkville", "success", "customer-story", "2026-03-02T10:03:23.971Z", "Parkville Creekside Development Relies on AT&T Business Fiber", "https://www.business.att.com/learn/customer-stories/parkville.html")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}
